# Apply the "merged to master" edits to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: TestStepExecution
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# New Timeout value for the 4th click step.
$ws1.Range("G4").Value = 2

# Replace hard-coded test data with placeholder tokens.
$ws1.Range("D7").Value = '$(Email)'
$ws1.Range("D9").Value = '$(Name)'
$ws1.Range("D11").Value = '$(Message)'

# Turn the A1:I13 block into a real Excel table ("Tabelle1").
$tbl1 = $ws1.ListObjects.Add(1, $ws1.Range("A1:I13"), $null, 1)
$tbl1.Name = "Tabelle1"
$tbl1.TableStyle = "TableStyleMedium1"

# Mirror the header-row shading/bold that Excel applies with the table style.
$hdr1 = $ws1.Range("A1:I1")
$hdrFont1 = $hdr1.Font
$hdrFont1.Bold = $true
$hdrFont1.Size = 9
$hdr1.Interior.Color = 5855577

# Column widths tweaked alongside the table.
$ws1.Columns.Item(2).ColumnWidth = 10.5
$ws1.Columns.Item(3).ColumnWidth = 33.5
$ws1.Columns.Item(4).ColumnWidth = 24.1
$ws1.Columns.Item(5).ColumnWidth = 10.1

# ---------------------------------------------------------------------
# Sheet 2: data
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Drop the old "Dummy" filler rows (8-11) before rebuilding the sheet.
$ws2.Range("A8:A11").ClearContents()

# Header row.
$ws2.Range("A1").Value = "Email"
$ws2.Range("B1").Value = "Name"
$ws2.Range("C1").Value = "Message"

# Data rows.
$emails = @("Good@baangt.org", "Tests@baangt.org", "Let@baangt.org", "you@baangt.org", "sleep@baangt.org", "well@baangt.org")
$names = @("Good", "Tests", "Let", "you", "Sleep", "Well")
$message = "Good Tests let you sleep well"

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $cellA = $ws2.Range("A$r")
    $cellA.Value = $emails[$i]
    $cellA.Style = "Link"
    $ws2.Hyperlinks.Add($cellA, "mailto:" + $emails[$i])
    $ws2.Range("B$r").Value = $names[$i]
    $ws2.Range("C$r").Value = $message
}

# Column widths for the new layout.
$ws2.Columns.Item(1).ColumnWidth = 15.1
$ws2.Columns.Item(3).ColumnWidth = 26.6

# ---------------------------------------------------------------------
# Selection / active-cell bookkeeping (matches the saved view state).
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C8").Select()
$ws1.Activate()
$ws1.Range("G5").Select()
